$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in column B (rows 2-7)
$ws.Range("B2").Value = 0.0230276641589093
$ws.Range("B3").Value = 0.0202303338201754
$ws.Range("B4").Value = 0.0000634954274860123
$ws.Range("B5").Value = 0.11609268046095
$ws.Range("B6").Value = 0.0175501124960426
$ws.Range("B7").Value = 0.0158290896823424

# Adjust column B width (target stored width ~21.19 chars; Excel quantizes
# ColumnWidth to whole-pixel increments, so 20.33 is the input that lands on
# the closest reachable stored width, 21.1667)
$ws.Columns.Item(2).ColumnWidth = 20.33

# Update selection to D10
$ws.Range("D10").Select()
